# trellis_statistics.xlsx - "Add files via upload"
# Adds a new "sum" style summary row-set to cluster30_BB (3 rows covering the
# MF/20BB/10BB totals) and refreshes the existing "sum" rows on cluster20_BB
# and cluster10_BB with updated totals, then leaves the workbook's view state
# (active sheet/zoom/selection) the way the author left it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# cluster30_BB: append three new "sum" rows (34-36) below the existing sum
# row (33), re-using that row's formatting.
# ---------------------------------------------------------------------------
$wsBB30 = $wb.Worksheets.Item("cluster30_BB")

$wsBB30.Range("A33:F33").Copy() | Out-Null
$wsBB30.Range("A34:F36").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 34 - matches the refreshed cluster20_BB sum totals
$wsBB30.Range("A34").Value = "sum"
$wsBB30.Range("B34").Value = "694/2664/1557"
$wsBB30.Range("C34").Value = "189/667/444"
$wsBB30.Range("E34").Value = "734/2664/2122"
$wsBB30.Range("F34").Value = "185/667/461"

# Row 35 - matches the refreshed cluster10_BB sum totals
$wsBB30.Range("A35").Value = "sum"
$wsBB30.Range("B35").Value = "697/2664/1558"
$wsBB30.Range("C35").Value = "184/667/425"
$wsBB30.Range("E35").Value = "726/2664/2036"
$wsBB30.Range("F35").Value = "182/667/476"

# Row 36 - matches cluster30_MF's own sum totals (same figures, reordered to
# this sheet's B/C/E/F column layout). Its F cell uses the plain (non-text)
# style like column C, not the s="8" text style the other two rows use, so
# copy that formatting across after the bulk paste above.
$wsBB30.Range("A36").Value = "sum"
$wsBB30.Range("B36").Value = "679/2664/1487"
$wsBB30.Range("C36").Value = "178/667/407"
$wsBB30.Range("E36").Value = "709/2664/2087"
$wsBB30.Range("F36").Value = "177/667/453"

$wsBB30.Range("C36").Copy() | Out-Null
$wsBB30.Range("F36").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$wsBB30.Range("F36").Value = "177/667/453"

# ---------------------------------------------------------------------------
# cluster20_BB: refresh the existing sum row (23) totals in place.
# ---------------------------------------------------------------------------
$wsBB20 = $wb.Worksheets.Item("cluster20_BB")
$wsBB20.Range("B23").Value = "694/2664/1557"
$wsBB20.Range("C23").Value = "189/667/444"
$wsBB20.Range("E23").Value = "734/2664/2122"
$wsBB20.Range("F23").Value = "185/667/461"

# ---------------------------------------------------------------------------
# cluster10_BB: refresh the existing sum row (13) totals in place.
# ---------------------------------------------------------------------------
$wsBB10 = $wb.Worksheets.Item("cluster10_BB")
$wsBB10.Range("B13").Value = "697/2664/1558"
$wsBB10.Range("C13").Value = "184/667/425"
$wsBB10.Range("E13").Value = "726/2664/2036"
$wsBB10.Range("F13").Value = "182/667/476"

# ---------------------------------------------------------------------------
# View state: restore each sheet's zoom/selection, and leave cluster30_BB as
# the active tab (matches the saved workbookView/activeTab + per-sheet
# tabSelected flag).
# ---------------------------------------------------------------------------
$wsMF30 = $wb.Worksheets.Item("cluster30_MF")
$wsMF30.Activate()
$excel.ActiveWindow.Zoom = 85
$wsMF30.Range("A32:F32").Select() | Out-Null

$wsBB20.Activate()
$wsBB20.Range("E23").Select() | Out-Null

$wsBB10.Activate()
$wsBB10.Range("E17").Select() | Out-Null

$wsBB30.Activate()
$excel.ActiveWindow.Zoom = 115
$wsBB30.Range("F33").Select() | Out-Null
